$d = $word.ActiveDocument

# Locate the "SourceCode" styled paragraph holding the stray R console
# output ("## null device" / "##  1") that follows the Results
# placeholder text, and remove the whole paragraph (including its
# trailing paragraph mark) now that the methodology/plot is still being
# drafted.
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    if ($text -like "*null device*") {
        $p.Range.Delete()
        break
    }
}
